# Append three new rows of raw voltage-read data to the "Voltage Results"
# sheet (rows 22-24), matching the lab's newer export format (date/time
# column now uses the m/d/yy h:mm display format instead of the old
# m/d/yyyy hh:mm:ss one), then leave the selection where the operator's
# cursor ended up after typing the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Voltage Results")

# New raw data rows pulled from the instrument log.
$newRows = @(
    @{ Row = 22; Date = 42998.429166666669; Volt = 1.63;               Err = 0.002 },
    @{ Row = 23; Date = 43000.432523148149; Volt = 1.6240000000000001; Err = 0.002 },
    @{ Row = 24; Date = 43000.68472222222;  Volt = 0.90400000000000003;Err = 0.002 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value2 = $r.Date
    $ws.Cells.Item($rowNum, 1).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($rowNum, 2).Value2 = $r.Volt
    $ws.Cells.Item($rowNum, 3).Value2 = $r.Err
}

# The file was relocated out of the synced Dropbox folder into a local,
# version-controlled "Documents\ResearchCode" tree (see commit message -
# raw data files are being pulled out of cloud sync to save space).
try {
    $wb.SaveAs("C:\Users\Michael\Documents\ResearchCode\Redoxotron\Data\Cr-Hematite_C-Ferrocyanide Redoxotron Results.xlsx")
} catch {
}

# Leave the cursor on B26, matching where the operator's selection ended
# up after entering the new rows.
$ws.Range("B26").Select() | Out-Null
